$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 'D2' '247.34'
Set-TextCell 'E2' '1.36%'
Set-TextCell 'D3' '30.14'
Set-TextCell 'E3' '11.21%'
Set-TextCell 'D4' '5.179'
Set-TextCell 'E4' '0.45%'
Set-TextCell 'D5' '0.05742'
Set-TextCell 'E5' '2.20%'
Set-TextCell 'D6' '6.592'
Set-TextCell 'E6' '1.72%'
Set-TextCell 'D7' '0.8573'
Set-TextCell 'E7' '5.04%'
Set-TextCell 'D8' '0.8720'
Set-TextCell 'E8' '4.76%'
Set-TextCell 'D9' '0.1365'
Set-TextCell 'E9' '2.85%'
Set-TextCell 'D10' '0.07067'
Set-TextCell 'E10' '2.40%'
Set-TextCell 'D11' '0.02923'
Set-TextCell 'E11' '1.21%'
Set-TextCell 'D12' '0.09389'
Set-TextCell 'E12' '0.06%'
Set-TextCell 'D13' '0.001515'
Set-TextCell 'E13' '-0.79%'
Set-TextCell 'D14' '0.04139'
Set-TextCell 'E14' '-7.89%'
Set-TextCell 'D15' '0.0006025'
Set-TextCell 'E15' '0.96%'
Set-TextCell 'D16' '0.006174'
Set-TextCell 'E16' '0.06%'
Set-TextCell 'D17' '3.505'
Set-TextCell 'E17' '-2.90%'
Set-TextCell 'D18' '3.041'
Set-TextCell 'E18' '0.63%'
Set-TextCell 'D19' '2.181'
Set-TextCell 'E19' '-2.07%'
Set-TextCell 'E20' '2.33%'
Set-TextCell 'D21' '0.03288'
Set-TextCell 'E21' '6.47%'
Set-TextCell 'D22' '0.1307'
Set-TextCell 'E22' '1.25%'
Set-TextCell 'D23' '3.630'
Set-TextCell 'E23' '-2.96%'
Set-TextCell 'D24' '0.1379'
Set-TextCell 'E24' '0.39%'
Set-TextCell 'E25' '-0.91%'
Set-TextCell 'D26' '0.004504'
Set-TextCell 'E26' '0.35%'
Set-TextCell 'E27' '20.38%'
Set-TextCell 'D28' '0.0001390'
Set-TextCell 'E28' '-0.62%'
Set-TextCell 'D40' '0.03785'
Set-TextCell 'E40' '4.14%'
Set-TextCell 'B41' 'BKEXToken'
Set-TextCell 'C41' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell 'D41' '0.1073'
Set-TextCell 'E41' '-21.04%'
Set-TextCell 'B42' 'CEJI'
Set-TextCell 'C42' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell 'D42' '0.002588'
Set-TextCell 'E42' '-1.55%'
Set-TextCell 'B43' 'KickToken'
Set-TextCell 'C43' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell 'D43' '0.003527'
Set-TextCell 'E43' '-41.77%'
Set-TextCell 'D44' '0.009586'
Set-TextCell 'E44' '16.98%'
Set-TextCell 'D45' '0.00005099'
Set-TextCell 'E45' '-4.08%'
Set-TextCell 'E46' '-0.05%'
Set-TextCell 'D47' '0.08893'
Set-TextCell 'E47' '-18.39%'
Set-TextCell 'D48' '0.002744'
Set-TextCell 'E48' '10.16%'
Set-TextCell 'E49' '-0.05%'
Set-TextCell 'E50' '-0.05%'
